$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-135 down to 20-136
$ws.Rows("19:19").Insert()

# Populate the new row 19 with the new record's data
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "Terminal La Palmera de La Serena"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = "2022-11-16"
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112052
$ws.Range("G19").Value = "Albahaca"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 1200
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 4500
$ws.Range("M19").Value = 4250
$ws.Range("N19").Value = "$/paquete"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 4250
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"
